$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures (includes the Cronos / EnergySwap
# row-order swap at rows 49-50). Values are forced to Text so Excel does not
# reinterpret numeric-looking strings (e.g. "208.99") as numbers, matching the
# original inline-string cell contents.
$updates = [ordered]@{
    "D2" = "25.989.24"
    "E2" = "  -0.56%  "
    "D3" = "1.640.89"
    "E3" = "  -1.74%  "
    "D4" = "1.002"
    "E4" = "  -0.08%  "
    "D5" = "208.99"
    "E5" = "  -0.95%  "
    "D6" = "0.5152"
    "E6" = "  -1.40%  "
    "D7" = "1.002"
    "E7" = "  -0.08%  "
    "D8" = "0.2560"
    "E8" = "  -3.11%  "
    "D9" = "0.06216"
    "E9" = "  -0.73%  "
    "D10" = "20.36"
    "E10" = "  -3.81%  "
    "D11" = "0.07542"
    "E11" = "  +0.33%  "
    "D12" = "1.641.69"
    "E12" = "  -1.63%  "
    "D13" = "4.355"
    "E13" = "  -1.94%  "
    "D14" = "1.862.83"
    "E14" = "  -1.85%  "
    "D15" = "0.5383"
    "E15" = "  -3.81%  "
    "D16" = "0.0₅7961"
    "E16" = "  -0.16%  "
    "D17" = "64.96"
    "E17" = "  -2.05%  "
    "D18" = "26.012.62"
    "E18" = "  -0.66%  "
    "E19" = "  -0.10%  "
    "D20" = "4.644"
    "E20" = "  -3.00%  "
    "D21" = "185.63"
    "E21" = "  -0.80%  "
    "D22" = "10.02"
    "E22" = "  -3.21%  "
    "D23" = "6.085"
    "E23" = "  -1.44%  "
    "E24" = "  -0.08%  "
    "D25" = "145.30"
    "E25" = "  -1.76%  "
    "E26" = "  -3.34%  "
    "D27" = "0.1190"
    "E27" = "  -4.41%  "
    "D28" = "15.44"
    "E28" = "  -3.03%  "
    "D29" = "1.374"
    "D30" = "0.05956"
    "E30" = "  -3.99%  "
    "E31" = "  -2.99%  "
    "D32" = "3.356"
    "E32" = "  -3.31%  "
    "D33" = "3.344"
    "E33" = "  -2.49%  "
    "D34" = "1.606"
    "E34" = "  -0.34%  "
    "D35" = "0.9682"
    "E35" = "  -2.40%  "
    "D36" = "2.373"
    "E36" = "  -1.31%  "
    "D37" = "2.726"
    "E37" = "  +0.65%  "
    "D38" = "0.5821"
    "E38" = "  -3.82%  "
    "D39" = "0.01586"
    "E39" = "  -1.60%  "
    "D40" = "1.047.92"
    "E40" = "  -1.97%  "
    "D41" = "5.766"
    "E41" = "  -5.94%  "
    "D42" = "1.002"
    "E42" = "  -0.28%  "
    "D43" = "0.8389"
    "E43" = "  -2.99%  "
    "D44" = "99.67"
    "E44" = "  +0.09%  "
    "D45" = "1.786.87"
    "E45" = "  -1.79%  "
    "E46" = "  -1.79%  "
    "E47" = "  +0.31%  "
    "D48" = "54.16"
    "E48" = "  -3.20%  "
    "B49" = "EnergySwap"
    "C49" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D49" = "7.929"
    "E49" = "  -0.04%  "
    "B50" = "Cronos"
    "C50" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D50" = "0.05197"
    "E50" = "  -0.97%  "
    "D51" = "0.4230"
    "E51" = "  -0.52%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
